$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.722.31"
$ws.Range("E2").Value = "  -1.77%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.384.94"
$ws.Range("E3").Value = "  -2.15%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.33"
$ws.Range("E5").Value = "  -2.22%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.86"
$ws.Range("E6").Value = "  -3.74%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.08%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.386.65"
$ws.Range("E8").Value = "  -2.12%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.27%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -2.28%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -2.14%  "

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.402"
$ws.Range("E12").Value = "  +2.58%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.964.74"
$ws.Range("E13").Value = "  -2.12%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.43"
$ws.Range("E14").Value = "  +1.51%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +1.56%  "

# Row 16 - was ShibaInu, now WrappedEther
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.406.13"
$ws.Range("E16").Value = "  -1.49%  "

# Row 17 - was WrappedEther, now ShibaInu
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000171"
$ws.Range("E17").Value = "  -2.30%  "

# Row 18 - WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.838.70"
$ws.Range("E18").Value = "  -1.72%  "

# Row 19 - Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.27"
$ws.Range("E19").Value = "  +0.32%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -2.19%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -5.88%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.14"
$ws.Range("E22").Value = "  -1.16%  "

# Row 23 - Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.563"
$ws.Range("E23").Value = "  -0.80%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.72"
$ws.Range("E24").Value = "  +0.11%  "

# Row 25 - Dai
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -5.60%  "

# Row 27 - WrappedeETH
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.522.71"
$ws.Range("E27").Value = "  -2.18%  "

# Row 28 - Kaspa
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.178"
$ws.Range("E28").Value = "  -2.84%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.66%  "

# Row 30 - RenderToken
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.45"
$ws.Range("E30").Value = "  -3.34%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.04"
$ws.Range("E31").Value = "  -2.14%  "

# Row 32 - was Fetch.AI, now PancakeSwap
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.15"
$ws.Range("E32").Value = "  -2.03%  "

# Row 33 - was PancakeSwap, now Fetch.AI
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.44"
$ws.Range("E33").Value = "  -3.25%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  -0.01%  "

# Row 35 - EthereumClassic
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.67"
$ws.Range("E35").Value = "  -2.36%  "

# Row 36 - Aptos
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.02"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37 - Monero
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.69"
$ws.Range("E37").Value = "  -0.21%  "

# Row 38 - NEARProtocol
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.03"
$ws.Range("E38").Value = "  -2.63%  "

# Row 39 - RenzoRestakedETH
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.418.81"
$ws.Range("E39").Value = "  -1.97%  "

# Row 40 - ImmutableX
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.50"
$ws.Range("E40").Value = "  -4.79%  "

# Row 41 - was EnergySwap, now Hedera
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0779"
$ws.Range("E41").Value = "  -0.69%  "

# Row 42 - was Hedera, now EnergySwap
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.83"
$ws.Range("E42").Value = "  +1.68%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  -2.95%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  +0.05%  "

# Row 45 - Filecoin
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.45"
$ws.Range("E45").Value = "  -1.58%  "

# Row 46 - OKB
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.75"
$ws.Range("E46").Value = "  -2.11%  "

# Row 47 - Stacks
$ws.Range("E47").Value = "  -2.96%  "

# Row 48 - Maker
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.531.81"
$ws.Range("E48").Value = "  -1.59%  "

# Row 49 - ONDO
$ws.Range("E49").Value = "  -4.50%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.61"
$ws.Range("E50").Value = "  +2.06%  "

# Row 51 - Cosmos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.88"
$ws.Range("E51").Value = "  -0.98%  "
